$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4664
$ws.Range("L3").Value = 5028
$ws.Range("H4").Value = 1759
$ws.Range("J4").Value = 1872
$ws.Range("K4").Value = 1782
$ws.Range("L4").Value = 1247
$ws.Range("L5").Value = 296
$ws.Range("L6").Value = 4283
$ws.Range("H7").Value = 26075
$ws.Range("J7").Value = 29348
$ws.Range("K7").Value = 27574
$ws.Range("L7").Value = 15518

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 297
$ws.Range("L7").Value = 1028

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 341

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 193
$ws.Range("L3").Value = 239
$ws.Range("L7").Value = 707

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 199
$ws.Range("L6").Value = 156
$ws.Range("L7").Value = 572

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 301

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 73
$ws.Range("L3").Value = 109
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L8").Value = 1028
$ws.Range("L10").Value = 102
$ws.Range("L11").Value = 250
$ws.Range("L12").Value = 36
$ws.Range("L14").Value = 82
$ws.Range("L17").Value = 29
$ws.Range("L19").Value = 427
$ws.Range("L20").Value = 392
$ws.Range("L22").Value = 45
$ws.Range("L23").Value = 169
$ws.Range("L25").Value = 91
$ws.Range("J29").Value = 1551
$ws.Range("L29").Value = 857
$ws.Range("L30").Value = 72
$ws.Range("L33").Value = 707
$ws.Range("L37").Value = 572
$ws.Range("L42").Value = 506
$ws.Range("L43").Value = 113
$ws.Range("L48").Value = 202
$ws.Range("L49").Value = 81
$ws.Range("L51").Value = 191
$ws.Range("L54").Value = 324
$ws.Range("L55").Value = 147
$ws.Range("L60").Value = 98
$ws.Range("H63").Value = 310
$ws.Range("K63").Value = 170
$ws.Range("L63").Value = 47
$ws.Range("L64").Value = 107
$ws.Range("L65").Value = 301
$ws.Range("L66").Value = 37
$ws.Range("L67").Value = 534
$ws.Range("L77").Value = 104
$ws.Range("L78").Value = 208
$ws.Range("L79").Value = 411
$ws.Range("L82").Value = 24
$ws.Range("L83").Value = 341
$ws.Range("L84").Value = 150
$ws.Range("L85").Value = 796
$ws.Range("L88").Value = 168
$ws.Range("L89").Value = 221
$ws.Range("L95").Value = 210
$ws.Range("L96").Value = 172
$ws.Range("L99").Value = 267
$ws.Range("H101").Value = 26075
$ws.Range("J101").Value = 29348
$ws.Range("K101").Value = 27574
$ws.Range("L101").Value = 15518

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 208
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 534

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 79
$ws.Range("L7").Value = 324

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 323
$ws.Range("J4").Value = 83
$ws.Range("L6").Value = 223
$ws.Range("J7").Value = 1551
$ws.Range("L7").Value = 857

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L4").Value = 41
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 133
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 427

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 36
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 170
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 506

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 56
$ws.Range("L3").Value = 67
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 42
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 172

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L4").Value = 29
$ws.Range("L7").Value = 411

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 118
$ws.Range("L7").Value = 392

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L4").Value = 20
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 46
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L4").Value = 18
$ws.Range("L6").Value = 39

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 54
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 240
$ws.Range("L3").Value = 323
$ws.Range("L5").Value = 19
$ws.Range("L7").Value = 796

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 36
